# Apply Ragnarok_Profits.xlsx market-price data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1642.6666
$ws.Range("I2").Value = 795
$ws.Range("K2").Value = 795
$ws.Range("M2").Value = -682
$ws.Range("H9").Value = 271.0909
$ws.Range("I9").Value = 176.14285
$ws.Range("K9").Value = 176.14285
$ws.Range("M9").Value = -7.14285000000001
$ws.Range("H40").Value = 55571908
$ws.Range("I40").Value = 1356.8
$ws.Range("J40").Value = 125035100
$ws.Range("K40").Value = 1356.8
$ws.Range("L40").Value = 125035100
$ws.Range("M40").Value = -1181.8
$ws.Range("N40").Value = -125035450
$ws.Range("H53").Value = 1313.6
$ws.Range("I53").Value = 1019.5714
$ws.Range("J53").Value = 1999.6666
$ws.Range("K53").Value = 1019.5714
$ws.Range("L53").Value = 1999.6666
$ws.Range("M53").Value = -382.5714
$ws.Range("N53").Value = -3273.6666
$ws.Range("H62").Value = 7363.75
$ws.Range("I62").Value = 6098.857
$ws.Range("K62").Value = 6098.857
$ws.Range("M62").Value = -5474.857
$ws.Range("H65").Value = 7363.75
$ws.Range("I65").Value = 6098.857
$ws.Range("K65").Value = 30494.285
$ws.Range("M65").Value = -27374.285
$ws.Range("H70").Value = 4066973.8
$ws.Range("J70").Value = 2900
$ws.Range("L70").Value = 8700
$ws.Range("N70").Value = -9240
$ws.Range("H73").Value = 4066973.8
$ws.Range("J73").Value = 2900
$ws.Range("L73").Value = 8700
$ws.Range("N73").Value = -10572
$ws.Range("H76").Value = 7777
$ws.Range("I76").Value = 8729
$ws.Range("K76").Value = 8729
$ws.Range("M76").Value = -8414
$ws.Range("H79").Value = 7777
$ws.Range("I79").Value = 8729
$ws.Range("K79").Value = 8729
$ws.Range("M79").Value = -7637
$ws.Range("H94").Value = 3067.7
$ws.Range("I94").Value = 3067.7
$ws.Range("K94").Value = 3067.7
$ws.Range("M94").Value = -2616.7
$ws.Range("H97").Value = 3305.2856
$ws.Range("J97").Value = 3022.8333
$ws.Range("L97").Value = 9068.499899999999
$ws.Range("N97").Value = -10060.4999
$ws.Range("H98").Value = 10417605
$ws.Range("I98").Value = 13889474
$ws.Range("K98").Value = 13889474
$ws.Range("M98").Value = -13887976
$ws.Range("H99").Value = 7686.1113
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 8396.875
$ws.Range("K99").Value = 6000
$ws.Range("L99").Value = 25190.625
$ws.Range("M99").Value = -4502
$ws.Range("N99").Value = -28186.625
$ws.Range("H103").Value = 38462852
$ws.Range("I103").Value = 821.1429000000001
$ws.Range("J103").Value = 83335224
$ws.Range("K103").Value = 2463.4287
$ws.Range("L103").Value = 250005672
$ws.Range("M103").Value = -1877.4287
$ws.Range("N103").Value = -250006844
$ws.Range("H106").Value = 6274.7617
$ws.Range("I106").Value = 6324.737
$ws.Range("J106").Value = 5800
$ws.Range("K106").Value = 6324.737
$ws.Range("L106").Value = 5800
$ws.Range("M106").Value = -5693.737
$ws.Range("N106").Value = -7062
$ws.Range("H122").Value = 10417605
$ws.Range("I122").Value = 13889474
$ws.Range("K122").Value = 41668422
$ws.Range("M122").Value = -41665972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 800.2
$ws.Range("I5").Value = 750.25
$ws.Range("K5").Value = 750.25
$ws.Range("M5").Value = -638.25
$ws.Range("H45").Value = 1559.1333
$ws.Range("I45").Value = 872.75
$ws.Range("J45").Value = 4304.6665
$ws.Range("K45").Value = 872.75
$ws.Range("L45").Value = 4304.6665
$ws.Range("M45").Value = -495.75
$ws.Range("N45").Value = -5058.6665
$ws.Range("H110").Value = 7145.0625
$ws.Range("I110").Value = 7375
$ws.Range("J110").Value = 6639.2
$ws.Range("K110").Value = 7375
$ws.Range("L110").Value = 6639.2
$ws.Range("M110").Value = -5330
$ws.Range("N110").Value = -10729.2
$ws.Range("H132").Value = 6252695
$ws.Range("I132").Value = 2716.7693
$ws.Range("J132").Value = 33335934
$ws.Range("K132").Value = 8150.3079
$ws.Range("L132").Value = 100007802
$ws.Range("M132").Value = -5620.3079
$ws.Range("N132").Value = -100012862

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 800.2
$ws.Range("I4").Value = 750.25
$ws.Range("K4").Value = 750.25
$ws.Range("M4").Value = -635.25
$ws.Range("H20").Value = 5222.9443
$ws.Range("I20").Value = 5365.6665
$ws.Range("J20").Value = 5080.222
$ws.Range("K20").Value = 5365.6665
$ws.Range("L20").Value = 5080.222
$ws.Range("M20").Value = -5118.6665
$ws.Range("N20").Value = -5574.222
$ws.Range("H86").Value = 2473.7144
$ws.Range("I86").Value = 1695.9375
$ws.Range("K86").Value = 1695.9375
$ws.Range("M86").Value = -572.9375
$ws.Range("H89").Value = 2473.7144
$ws.Range("I89").Value = 1695.9375
$ws.Range("K89").Value = 8479.6875
$ws.Range("M89").Value = -2863.6875
$ws.Range("H94").Value = 1937.84
$ws.Range("J94").Value = 1298.5
$ws.Range("L94").Value = 1298.5
$ws.Range("N94").Value = -2200.5
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 22232636
$ws.Range("I62").Value = 5866.1665
$ws.Range("J62").Value = 37050484
$ws.Range("K62").Value = 5866.1665
$ws.Range("L62").Value = 37050484
$ws.Range("M62").Value = -5242.1665
$ws.Range("N62").Value = -37051732
$ws.Range("H65").Value = 22232636
$ws.Range("I65").Value = 5866.1665
$ws.Range("J65").Value = 37050484
$ws.Range("K65").Value = 29330.8325
$ws.Range("L65").Value = 185252420
$ws.Range("M65").Value = -26210.8325
$ws.Range("N65").Value = -185258660
$ws.Range("H122").Value = 2428.5
$ws.Range("I122").Value = 2428.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7285.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -4835.5
$ws.Range("H132").Value = 2325.182
$ws.Range("I132").Value = 2325.182
$ws.Range("K132").Value = 6975.545999999999
$ws.Range("M132").Value = -4445.545999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 15084.875
$ws.Range("J69").Value = 27494.75
$ws.Range("L69").Value = 82484.25
$ws.Range("N69").Value = -84106.25
$ws.Range("H72").Value = 15084.875
$ws.Range("J72").Value = 27494.75
$ws.Range("L72").Value = 247452.75
$ws.Range("N72").Value = -255564.75
$ws.Range("H131").Value = 4603.92
$ws.Range("I131").Value = 1812.6364
$ws.Range("J131").Value = 6797.0713
$ws.Range("K131").Value = 5437.9092
$ws.Range("L131").Value = 20391.2139
$ws.Range("M131").Value = -397.9092000000001
$ws.Range("N131").Value = -30471.2139
$ws.Range("H141").Value = 5026.3335
$ws.Range("I141").Value = 5026.3335
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 15079.0005
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -9899.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2318986
$ws.Range("I113").Value = 3444.5
$ws.Range("J113").Value = 4634527
$ws.Range("K113").Value = 3444.5
$ws.Range("L113").Value = 4634527
$ws.Range("M113").Value = -1274.5
$ws.Range("N113").Value = -4638867
$ws.Range("H126").Value = 5357.077
$ws.Range("I126").Value = 4906.4
$ws.Range("J126").Value = 6859.3335
$ws.Range("K126").Value = 14719.2
$ws.Range("L126").Value = 20578.0005
$ws.Range("M126").Value = -12249.2
$ws.Range("N126").Value = -25518.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3162.125
$ws.Range("J22").Value = 4519.6
$ws.Range("L22").Value = 4519.6
$ws.Range("N22").Value = -5109.6
$ws.Range("H27").Value = 3162.125
$ws.Range("J27").Value = 4519.6
$ws.Range("L27").Value = 4519.6
$ws.Range("N27").Value = -4733.6
$ws.Range("H55").Value = 757.86664
$ws.Range("I55").Value = 294.7857
$ws.Range("J55").Value = 1163.0625
$ws.Range("K55").Value = 294.7857
$ws.Range("L55").Value = 1163.0625
$ws.Range("M55").Value = -121.7857
$ws.Range("N55").Value = -1509.0625
$ws.Range("H57").Value = 26402.7
$ws.Range("I57").Value = 26402.7
$ws.Range("K57").Value = 26402.7
$ws.Range("M57").Value = -25836.7
$ws.Range("H61").Value = 2827
$ws.Range("I61").Value = 2629.5454
$ws.Range("J61").Value = 4999
$ws.Range("K61").Value = 2629.5454
$ws.Range("L61").Value = 4999
$ws.Range("M61").Value = -2427.5454
$ws.Range("N61").Value = -5403
$ws.Range("H68").Value = 2318192.2
$ws.Range("J68").Value = 5836.875
$ws.Range("L68").Value = 5836.875
$ws.Range("N68").Value = -7334.875
$ws.Range("H71").Value = 2318192.2
$ws.Range("J71").Value = 5836.875
$ws.Range("L71").Value = 29184.375
$ws.Range("N71").Value = -36672.375
$ws.Range("H113").Value = 2827
$ws.Range("I113").Value = 2629.5454
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 2629.5454
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = -459.5454
$ws.Range("N113").Value = -9339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 34941
$ws.Range("J44").Value = 34941
$ws.Range("L44").Value = 34941
$ws.Range("N44").Value = -36049
$ws.Range("H49").Value = 34962
$ws.Range("J49").Value = 34962
$ws.Range("L49").Value = 34962
$ws.Range("N49").Value = -35422
$ws.Range("H140").Value = 65142.668
$ws.Range("J140").Value = 65142.668
$ws.Range("L140").Value = 65142.668
$ws.Range("N140").Value = -75502.66800000001
$ws.Range("H141").Value = 87814
$ws.Range("J141").Value = 87814
$ws.Range("L141").Value = 87814
$ws.Range("N141").Value = -98174
